$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5910.3335
$ws.Range("I76").Value = 4647.909
$ws.Range("J76").Value = 7299
$ws.Range("K76").Value = 4647.909
$ws.Range("L76").Value = 7299
$ws.Range("M76").Value = -4332.909
$ws.Range("N76").Value = -7929
$ws.Range("H79").Value = 5910.3335
$ws.Range("I79").Value = 4647.909
$ws.Range("J79").Value = 7299
$ws.Range("K79").Value = 4647.909
$ws.Range("L79").Value = 7299
$ws.Range("M79").Value = -3555.909
$ws.Range("N79").Value = -9483
$ws.Range("H116").Value = 7444.1665
$ws.Range("J116").Value = 8099.5713
$ws.Range("L116").Value = 8099.5713
$ws.Range("N116").Value = -14983.5713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 5114
$ws.Range("I28").Value = 5114
$ws.Range("K28").Value = 5114
$ws.Range("M28").Value = -4922
$ws.Range("H31").Value = 2900
$ws.Range("I31").Value = 2900
$ws.Range("K31").Value = 2900
$ws.Range("M31").Value = -2606
$ws.Range("H52").Value = 33304.5
$ws.Range("J52").Value = 55900
$ws.Range("L52").Value = 55900
$ws.Range("N52").Value = -56536
$ws.Range("H61").Value = 1829
$ws.Range("I61").Value = 1829
$ws.Range("K61").Value = 1829
$ws.Range("M61").Value = -1617
$ws.Range("H99").Value = 5114
$ws.Range("I99").Value = 5114
$ws.Range("K99").Value = 5114
$ws.Range("M99").Value = -2119
$ws.Range("H102").Value = 748.2308
$ws.Range("I102").Value = 384.27274
$ws.Range("K102").Value = 384.27274
$ws.Range("M102").Value = 1237.72726
$ws.Range("H136").Value = 1829
$ws.Range("I136").Value = 1829
$ws.Range("K136").Value = 5487
$ws.Range("M136").Value = -2937

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 299984.5
$ws.Range("I122").Value = 299979
$ws.Range("K122").Value = 299979
$ws.Range("M122").Value = -295079
$ws.Range("H134").Value = 2645.0625
$ws.Range("I134").Value = 2370.923
$ws.Range("J134").Value = 3833
$ws.Range("K134").Value = 7112.768999999999
$ws.Range("L134").Value = 11499
$ws.Range("M134").Value = -4577.768999999999
$ws.Range("N134").Value = -16569

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3169.1
$ws.Range("I58").Value = 2045.4286
$ws.Range("J58").Value = 3774.1538
$ws.Range("K58").Value = 2045.4286
$ws.Range("L58").Value = 3774.1538
$ws.Range("M58").Value = -1842.4286
$ws.Range("N58").Value = -4180.1538
$ws.Range("H122").Value = 3577.697
$ws.Range("I122").Value = 3759.0476
$ws.Range("J122").Value = 3260.3333
$ws.Range("K122").Value = 11277.1428
$ws.Range("L122").Value = 9780.999899999999
$ws.Range("M122").Value = -8827.1428
$ws.Range("N122").Value = -14680.9999
$ws.Range("H134").Value = 4190.4165
$ws.Range("I134").Value = 3345.75
$ws.Range("J134").Value = 4612.75
$ws.Range("K134").Value = 10037.25
$ws.Range("L134").Value = 13838.25
$ws.Range("M134").Value = -7502.25
$ws.Range("N134").Value = -18908.25
$ws.Range("H136").Value = 3169.1
$ws.Range("I136").Value = 2045.4286
$ws.Range("J136").Value = 3774.1538
$ws.Range("K136").Value = 6136.2858
$ws.Range("L136").Value = 11322.4614
$ws.Range("M136").Value = -3586.2858
$ws.Range("N136").Value = -16422.4614
$ws.Range("H141").Value = 28234.059
$ws.Range("J141").Value = 28234.059
$ws.Range("L141").Value = 28234.059
$ws.Range("N141").Value = -38594.059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 10000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H99").Value = 15794.2
$ws.Range("I99").Value = 8157
$ws.Range("J99").Value = 27250
$ws.Range("K99").Value = 8157
$ws.Range("L99").Value = 27250
$ws.Range("M99").Value = -5911
$ws.Range("N99").Value = -31742
$ws.Range("H102").Value = 1687.0333
$ws.Range("I102").Value = 458.66666
$ws.Range("J102").Value = 4553.222
$ws.Range("K102").Value = 458.66666
$ws.Range("L102").Value = 4553.222
$ws.Range("M102").Value = 1163.33334
$ws.Range("N102").Value = -7797.222
$ws.Range("H132").Value = 2820.2942
$ws.Range("I132").Value = 1927.6923
$ws.Range("J132").Value = 5721.25
$ws.Range("K132").Value = 5783.0769
$ws.Range("L132").Value = 17163.75
$ws.Range("M132").Value = -3253.0769
$ws.Range("N132").Value = -22223.75
$ws.Range("H133").Value = 100001
$ws.Range("J133").Value = 100001
$ws.Range("L133").Value = 100001
$ws.Range("N133").Value = -110121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3096.4443
$ws.Range("I82").Value = 3109.7144
$ws.Range("J82").Value = 3050
$ws.Range("K82").Value = 3109.7144
$ws.Range("L82").Value = 3050
$ws.Range("M82").Value = -2748.7144
$ws.Range("N82").Value = -3772
$ws.Range("H85").Value = 3096.4443
$ws.Range("I85").Value = 3109.7144
$ws.Range("J85").Value = 3050
$ws.Range("K85").Value = 3109.7144
$ws.Range("L85").Value = 3050
$ws.Range("M85").Value = -1861.7144
$ws.Range("N85").Value = -5546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 100000
$ws.Range("J46").Value = 100000
$ws.Range("L46").Value = 100000
$ws.Range("N46").Value = -100462
$ws.Range("H96").Value = 1231.875
$ws.Range("I96").Value = 1170
$ws.Range("K96").Value = 1170
$ws.Range("M96").Value = 203
$ws.Range("H107").Value = 787.3333
$ws.Range("I107").Value = 437.2
$ws.Range("K107").Value = 1311.6
$ws.Range("M107").Value = 608.4000000000001
$ws.Range("H132").Value = 2257.125
$ws.Range("I132").Value = 1594.5
$ws.Range("J132").Value = 4245
$ws.Range("K132").Value = 4783.5
$ws.Range("L132").Value = 12735
$ws.Range("M132").Value = -2253.5
$ws.Range("N132").Value = -17795
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 300000
$ws.Range("N134").Value = -305070
